$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 496
$ws.Range("J32").Value = 545
$ws.Range("L32").Value = 545
$ws.Range("N32").Value = -1197
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H51").Value = 2114.8462
$ws.Range("I51").Value = 2050
$ws.Range("J51").Value = 2143.6667
$ws.Range("K51").Value = 2050
$ws.Range("L51").Value = 2143.6667
$ws.Range("M51").Value = -1566
$ws.Range("N51").Value = -3111.6667
$ws.Range("H129").Value = 1030.4744
$ws.Range("J129").Value = 1091.1538
$ws.Range("L129").Value = 3273.4614
$ws.Range("N129").Value = -13273.4614
$ws.Range("H132").Value = 1589.0588
$ws.Range("I132").Value = 1166.875
$ws.Range("J132").Value = 2602.3
$ws.Range("K132").Value = 3500.625
$ws.Range("L132").Value = 7806.900000000001
$ws.Range("M132").Value = -970.625
$ws.Range("N132").Value = -12866.9
$ws.Range("H135").Value = 1538.9692
$ws.Range("I135").Value = 1136.1731
$ws.Range("J135").Value = 3150.1538
$ws.Range("K135").Value = 10225.5579
$ws.Range("L135").Value = 28351.3842
$ws.Range("M135").Value = -7690.5579
$ws.Range("N135").Value = -33421.3842
$ws.Range("H137").Value = 1481.1372
$ws.Range("I137").Value = 1258.8485
$ws.Range("J137").Value = 1888.6666
$ws.Range("K137").Value = 3776.5455
$ws.Range("L137").Value = 5665.9998
$ws.Range("M137").Value = -1226.5455
$ws.Range("N137").Value = -10765.9998
$ws.Range("H138").Value = 3462.5393
$ws.Range("I138").Value = 1227.68
$ws.Range("J138").Value = 4335.5312
$ws.Range("K138").Value = 3683.04
$ws.Range("L138").Value = 13006.5936
$ws.Range("M138").Value = 1456.96
$ws.Range("N138").Value = -23286.5936
$ws.Range("H141").Value = 1802
$ws.Range("I141").Value = 1444.68
$ws.Range("J141").Value = 2397.5334
$ws.Range("K141").Value = 4334.04
$ws.Range("L141").Value = 7192.600199999999
$ws.Range("M141").Value = 845.96
$ws.Range("N141").Value = -17552.6002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5697.53
$ws.Range("I32").Value = 4172.0864
$ws.Range("J32").Value = 12200.737
$ws.Range("K32").Value = 4172.0864
$ws.Range("L32").Value = 12200.737
$ws.Range("M32").Value = -3885.0864
$ws.Range("N32").Value = -12774.737
$ws.Range("H74").Value = 1004.44446
$ws.Range("I74").Value = 845.8372000000001
$ws.Range("J74").Value = 1624.4546
$ws.Range("K74").Value = 845.8372000000001
$ws.Range("L74").Value = 1624.4546
$ws.Range("M74").Value = 28.16279999999995
$ws.Range("N74").Value = -3372.4546
$ws.Range("H77").Value = 1004.44446
$ws.Range("I77").Value = 845.8372000000001
$ws.Range("J77").Value = 1624.4546
$ws.Range("K77").Value = 4229.186000000001
$ws.Range("L77").Value = 8122.273
$ws.Range("M77").Value = 138.8139999999994
$ws.Range("N77").Value = -16858.273
$ws.Range("H80").Value = 33892.332
$ws.Range("I80").Value = 17827.5
$ws.Range("J80").Value = 41924.75
$ws.Range("K80").Value = 17827.5
$ws.Range("L80").Value = 41924.75
$ws.Range("M80").Value = -16829.5
$ws.Range("N80").Value = -43920.75
$ws.Range("H81").Value = 31100
$ws.Range("I81").Value = 20000
$ws.Range("J81").Value = 42200
$ws.Range("K81").Value = 20000
$ws.Range("L81").Value = 42200
$ws.Range("M81").Value = -19002
$ws.Range("N81").Value = -44196
$ws.Range("H83").Value = 33892.332
$ws.Range("I83").Value = 17827.5
$ws.Range("J83").Value = 41924.75
$ws.Range("K83").Value = 53482.5
$ws.Range("L83").Value = 125774.25
$ws.Range("M83").Value = -48490.5
$ws.Range("N83").Value = -135758.25
$ws.Range("H84").Value = 31100
$ws.Range("I84").Value = 20000
$ws.Range("J84").Value = 42200
$ws.Range("K84").Value = 60000
$ws.Range("L84").Value = 126600
$ws.Range("M84").Value = -55008
$ws.Range("N84").Value = -136584
$ws.Range("H86").Value = 35280
$ws.Range("I86").Value = 35280
$ws.Range("K86").Value = 35280
$ws.Range("M86").Value = -34094
$ws.Range("H89").Value = 35280
$ws.Range("I89").Value = 35280
$ws.Range("K89").Value = 105840
$ws.Range("M89").Value = -99912
$ws.Range("H122").Value = 3206207.8
$ws.Range("I122").Value = 4274491.5
$ws.Range("J122").Value = 1357
$ws.Range("K122").Value = 12823474.5
$ws.Range("L122").Value = 4071
$ws.Range("M122").Value = -12821024.5
$ws.Range("N122").Value = -8971

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1863.5454
$ws.Range("I86").Value = 1800
$ws.Range("J86").Value = 1939.8
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 1939.8
$ws.Range("M86").Value = -677
$ws.Range("N86").Value = -4185.8
$ws.Range("H89").Value = 1863.5454
$ws.Range("I89").Value = 1800
$ws.Range("J89").Value = 1939.8
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 9699
$ws.Range("M89").Value = -3384
$ws.Range("N89").Value = -20931

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1224.7847
$ws.Range("I58").Value = 750.1795
$ws.Range("J58").Value = 1936.6923
$ws.Range("K58").Value = 750.1795
$ws.Range("L58").Value = 1936.6923
$ws.Range("M58").Value = -547.1795
$ws.Range("N58").Value = -2342.6923
$ws.Range("H107").Value = 8547654
$ws.Range("I107").Value = 11111607
$ws.Range("K107").Value = 11111607
$ws.Range("M107").Value = -11109687
$ws.Range("H132").Value = 1947.9636
$ws.Range("I132").Value = 1541.381
$ws.Range("J132").Value = 3261.5386
$ws.Range("K132").Value = 4624.143
$ws.Range("L132").Value = 9784.6158
$ws.Range("M132").Value = -2094.143
$ws.Range("N132").Value = -14844.6158
$ws.Range("H134").Value = 1581.4445
$ws.Range("I134").Value = 1678.5098
$ws.Range("J134").Value = 1345.7142
$ws.Range("K134").Value = 5035.5294
$ws.Range("L134").Value = 4037.1426
$ws.Range("M134").Value = -2500.5294
$ws.Range("N134").Value = -9107.142599999999
$ws.Range("H136").Value = 1224.7847
$ws.Range("I136").Value = 750.1795
$ws.Range("J136").Value = 1936.6923
$ws.Range("K136").Value = 2250.5385
$ws.Range("L136").Value = 5810.0769
$ws.Range("M136").Value = 299.4615000000003
$ws.Range("N136").Value = -10910.0769

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2979.1667
$ws.Range("J54").Value = 2979.1667
$ws.Range("L54").Value = 8937.500100000001
$ws.Range("N54").Value = -10055.5001
$ws.Range("H131").Value = 2000934.2
$ws.Range("I131").Value = 16667017
$ws.Range("J131").Value = 1013.86365
$ws.Range("K131").Value = 50001051
$ws.Range("L131").Value = 3041.59095
$ws.Range("M131").Value = -49996011
$ws.Range("N131").Value = -13121.59095

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 36719600
$ws.Range("I122").Value = 66552480
$ws.Range("J122").Value = 2209.4614
$ws.Range("K122").Value = 199657440
$ws.Range("L122").Value = 6628.3842
$ws.Range("M122").Value = -199654990
$ws.Range("N122").Value = -11528.3842
$ws.Range("H123").Value = 22694.375
$ws.Range("J123").Value = 22694.375
$ws.Range("L123").Value = 22694.375
$ws.Range("N123").Value = -27594.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 840.8
$ws.Range("I16").Value = 840.8
$ws.Range("K16").Value = 840.8
$ws.Range("M16").Value = -670.8
$ws.Range("H68").Value = 50001836
$ws.Range("I68").Value = 1992.8572
$ws.Range("J68").Value = 166668130
$ws.Range("K68").Value = 1992.8572
$ws.Range("L68").Value = 166668130
$ws.Range("M68").Value = -1243.8572
$ws.Range("N68").Value = -166669628
$ws.Range("H71").Value = 50001836
$ws.Range("I71").Value = 1992.8572
$ws.Range("J71").Value = 166668130
$ws.Range("K71").Value = 9964.286
$ws.Range("L71").Value = 833340650
$ws.Range("M71").Value = -6220.286
$ws.Range("N71").Value = -833348138
$ws.Range("H82").Value = 1003668.2
$ws.Range("I82").Value = 1112372.2
$ws.Range("J82").Value = 514500
$ws.Range("K82").Value = 1112372.2
$ws.Range("L82").Value = 514500
$ws.Range("M82").Value = -1112011.2
$ws.Range("N82").Value = -515222
$ws.Range("H85").Value = 1003668.2
$ws.Range("I85").Value = 1112372.2
$ws.Range("J85").Value = 514500
$ws.Range("K85").Value = 1112372.2
$ws.Range("L85").Value = 514500
$ws.Range("M85").Value = -1111124.2
$ws.Range("N85").Value = -516996
$ws.Range("H122").Value = 3880789.2
$ws.Range("I122").Value = 4467592
$ws.Range("K122").Value = 13402776
$ws.Range("M122").Value = -13400326
$ws.Range("H128").Value = 31000
$ws.Range("J128").Value = 31000
$ws.Range("L128").Value = 31000
$ws.Range("N128").Value = -40960
$ws.Range("H132").Value = 12148452
$ws.Range("I132").Value = 14445954
$ws.Range("J132").Value = 4514.7144
$ws.Range("K132").Value = 43337862
$ws.Range("L132").Value = 13544.1432
$ws.Range("M132").Value = -43335332
$ws.Range("N132").Value = -18604.1432
$ws.Range("H136").Value = 6580.537
$ws.Range("I136").Value = 4276.9287
$ws.Range("J136").Value = 14643.167
$ws.Range("K136").Value = 12830.7861
$ws.Range("L136").Value = 43929.501
$ws.Range("M136").Value = -10280.7861
$ws.Range("N136").Value = -49029.501

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2643
$ws.Range("I122").Value = 2569.1304
$ws.Range("J122").Value = 2885.7144
$ws.Range("K122").Value = 7707.3912
$ws.Range("L122").Value = 8657.143199999999
$ws.Range("M122").Value = -5257.3912
$ws.Range("N122").Value = -13557.1432
$ws.Range("H132").Value = 16499.375
$ws.Range("I132").Value = 20252.883
$ws.Range("J132").Value = 1774.0769
$ws.Range("K132").Value = 60758.649
$ws.Range("L132").Value = 5322.2307
$ws.Range("M132").Value = -58228.649
$ws.Range("N132").Value = -10382.2307
$ws.Range("H136").Value = 9436891
$ws.Range("I136").Value = 3730.4644
$ws.Range("J136").Value = 20002032
$ws.Range("K136").Value = 11191.3932
$ws.Range("L136").Value = 60006096
$ws.Range("M136").Value = -8641.393199999999
$ws.Range("N136").Value = -60011196
